$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 to I1:J1 so the new headers match existing style
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I and J columns (rows 2-12)
$values = @(
    @(8, 8),
    @(8, 9),
    @(6, 7),
    @(11, 11),
    @(5, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(7, 7),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
